$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: "No Technique" stays the same, model and score change
$ws.Range("B2").Value = "Naive Bayes"
$ws.Range("C2").Value = "0.99 ± 0.01"

# Row 3: technique, model, and score change
$ws.Range("A3").Value = "No Technique"
$ws.Range("B3").Value = "SVM"
$ws.Range("C3").Value = "0.98 ± 0.02"
